# Generate Report for Handback
#
# The localization-status report is refreshed after a handback attempt for
# file "b70d7998-5e51-43ca-b5ba-31e62ba3ef83" failed (the returned handback
# file name didn't match the expected handoff file name). Update the status
# everywhere it is surfaced and record the error detail on each language
# sheet.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh       = $wb.Worksheets.Item("zh-cn")
$wsDe       = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

# Overview sheet: row 3 is the b70d7998-... file; its zh-cn / de-de status
# columns (E, F) flip from "Ready for handoff" to the failure status.
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# Each language sheet's own "Status" column (C) for that same file.
$wsZh.Range("C3").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

# Record why the handback failed in the "Error Detail" column (P) for that
# row on each language sheet.
$wsZh.Range("P3").Value = "Handback file name: 0ckp1hwl.soo is different with handoff file name: b70d7998-5e51-43ca-b5ba-31e62ba3ef83.6fc89e8d1ca47346879ff53f15b2320052d21b4f.zh-cn."
$wsDe.Range("P3").Value = "Handback file name: 0ckp1hwl.soo is different with handoff file name: b70d7998-5e51-43ca-b5ba-31e62ba3ef83.6fc89e8d1ca47346879ff53f15b2320052d21b4f.de-de."

# Widen the Error Detail column so the long message is readable
# (width 40, matching the other wide text columns on these sheets).
$wsZh.Columns.Item(16).ColumnWidth = 39.17
$wsDe.Columns.Item(16).ColumnWidth = 39.17
